$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 29.20950566666667
$ws.Cells.Item(2, 8).Value = 87.628517
$ws.Cells.Item(2, 9).Value = 0.01829497698069002
$ws.Cells.Item(2, 10).Value = 0.01840828041918582
$ws.Cells.Item(2, 13).Value = 2.641124333333333
$ws.Cells.Item(2, 14).Value = 7.923373
$ws.Cells.Item(2, 15).Value = 0.2973163964900306
$ws.Cells.Item(2, 16).Value = 0.3204448291201116
$ws.Cells.Item(2, 17).Value = 77.14593618087123
$ws.Cells.Item(2, 18).Value = 694.313425627841
$ws.Cells.Item(2, 19).Value = 0.005439396629766816
$ws.Cells.Item(2, 20).Value = 0.005898838273321096

$ws.Cells.Item(3, 7).Value = 29.20950566666667
$ws.Cells.Item(3, 8).Value = 87.628517
$ws.Cells.Item(3, 9).Value = 0.01829497698069002
$ws.Cells.Item(3, 10).Value = 0.01840828041918582
$ws.Cells.Item(3, 15).Value = 0.258118258924649
$ws.Cells.Item(3, 16).Value = 0.2781974433645584
$ws.Cells.Item(3, 17).Value = 66.97503052370767
$ws.Cells.Item(3, 18).Value = 602.7752747133691
$ws.Cells.Item(3, 19).Value = 0.00472226760532224
$ws.Cells.Item(3, 20).Value = 0.005121136549355357

$ws.Cells.Item(4, 7).Value = 29.20950566666667
$ws.Cells.Item(4, 8).Value = 87.628517
$ws.Cells.Item(4, 9).Value = 0.01829497698069002
$ws.Cells.Item(4, 10).Value = 0.01840828041918582
$ws.Cells.Item(4, 13).Value = 0.7354653333333333
$ws.Cells.Item(4, 14).Value = 2.206396
$ws.Cells.Item(4, 15).Value = 0.08279273334096697
$ws.Cells.Item(4, 16).Value = 0.08923323301721346
$ws.Cells.Item(4, 17).Value = 21.48257882163689
$ws.Cells.Item(4, 18).Value = 193.343209394732
$ws.Cells.Item(4, 19).Value = 0.001514691150641398
$ws.Cells.Item(4, 20).Value = 0.001642630376091416

$ws.Cells.Item(5, 7).Value = 29.20950566666667
$ws.Cells.Item(5, 8).Value = 87.628517
$ws.Cells.Item(5, 9).Value = 0.01829497698069002
$ws.Cells.Item(5, 10).Value = 0.01840828041918582
$ws.Cells.Item(5, 13).Value = 1.9234645
$ws.Cells.Item(5, 14).Value = 3.846929
$ws.Cells.Item(5, 15).Value = 0.2165280621964276
$ws.Cells.Item(5, 16).Value = 0.1555812790893729
$ws.Cells.Item(5, 17).Value = 56.18344721238217
$ws.Cells.Item(5, 18).Value = 337.100683274293
$ws.Cells.Item(5, 19).Value = 0.003961375913557059
$ws.Cells.Item(5, 20).Value = 0.002863983813452787

$ws.Cells.Item(6, 7).Value = 29.20950566666667
$ws.Cells.Item(6, 8).Value = 87.628517
$ws.Cells.Item(6, 9).Value = 0.01829497698069002
$ws.Cells.Item(6, 10).Value = 0.01840828041918582
$ws.Cells.Item(6, 13).Value = 1.290238
$ws.Cells.Item(6, 14).Value = 3.870714
$ws.Cells.Item(6, 15).Value = 0.1452445490479259
$ws.Cells.Item(6, 16).Value = 0.1565432154087437
$ws.Cells.Item(6, 17).Value = 37.68721417234867
$ws.Cells.Item(6, 18).Value = 339.184927551138
$ws.Cells.Item(6, 19).Value = 0.002657245681402508
$ws.Cells.Item(6, 20).Value = 0.002881691406965164

$ws.Cells.Item(7, 9).Value = 0.913374480506715
$ws.Cells.Item(7, 10).Value = 0.9190311407684336
$ws.Cells.Item(7, 13).Value = 2.641124333333333
$ws.Cells.Item(7, 14).Value = 7.923373
$ws.Cells.Item(7, 15).Value = 0.2973163964900306
$ws.Cells.Item(7, 16).Value = 0.3204448291201116
$ws.Cells.Item(7, 17).Value = 3851.501396081551
$ws.Cells.Item(7, 18).Value = 34663.51256473396
$ws.Cells.Item(7, 19).Value = 0.2715612091902102
$ws.Cells.Item(7, 20).Value = 0.2944987768596019

$ws.Cells.Item(8, 9).Value = 0.913374480506715
$ws.Cells.Item(8, 10).Value = 0.9190311407684336
$ws.Cells.Item(8, 15).Value = 0.258118258924649
$ws.Cells.Item(8, 16).Value = 0.2781974433645584
$ws.Cells.Item(8, 19).Value = 0.235758630654599
$ws.Cells.Item(8, 20).Value = 0.2556721137341919

$ws.Cells.Item(9, 9).Value = 0.913374480506715
$ws.Cells.Item(9, 10).Value = 0.9190311407684336
$ws.Cells.Item(9, 13).Value = 0.7354653333333333
$ws.Cells.Item(9, 14).Value = 2.206396
$ws.Cells.Item(9, 15).Value = 0.08279273334096697
$ws.Cells.Item(9, 16).Value = 0.08923323301721346
$ws.Cells.Item(9, 17).Value = 1072.515111216997
$ws.Cells.Item(9, 18).Value = 9652.636000952971
$ws.Cells.Item(9, 19).Value = 0.07562076980503669
$ws.Cells.Item(9, 20).Value = 0.08200811993426514

$ws.Cells.Item(10, 9).Value = 0.913374480506715
$ws.Cells.Item(10, 10).Value = 0.9190311407684336
$ws.Cells.Item(10, 13).Value = 1.9234645
$ws.Cells.Item(10, 14).Value = 3.846929
$ws.Cells.Item(10, 15).Value = 0.2165280621964276
$ws.Cells.Item(10, 16).Value = 0.1555812790893729
$ws.Cells.Item(10, 17).Value = 2804.951706954842
$ws.Cells.Item(10, 18).Value = 16829.71024172905
$ws.Cells.Item(10, 19).Value = 0.1977712063237877
$ws.Cells.Item(10, 20).Value = 0.1429840404037184

$ws.Cells.Item(11, 9).Value = 0.913374480506715
$ws.Cells.Item(11, 10).Value = 0.9190311407684336
$ws.Cells.Item(11, 13).Value = 1.290238
$ws.Cells.Item(11, 14).Value = 3.870714
$ws.Cells.Item(11, 15).Value = 0.1452445490479259
$ws.Cells.Item(11, 16).Value = 0.1565432154087437
$ws.Cells.Item(11, 17).Value = 1881.529542384588
$ws.Cells.Item(11, 18).Value = 16933.7658814613
$ws.Cells.Item(11, 19).Value = 0.1326626645330815
$ws.Cells.Item(11, 20).Value = 0.1438680898366563

$ws.Cells.Item(12, 7).Value = 57.98602933333333
$ws.Cells.Item(12, 8).Value = 173.958088
$ws.Cells.Item(12, 9).Value = 0.03631876156896331
$ws.Cells.Item(12, 10).Value = 0.03654368891224535
$ws.Cells.Item(12, 13).Value = 2.641124333333333
$ws.Cells.Item(12, 14).Value = 7.923373
$ws.Cells.Item(12, 15).Value = 0.2973163964900306
$ws.Cells.Item(12, 16).Value = 0.3204448291201116
$ws.Cells.Item(12, 17).Value = 153.1483130656471
$ws.Cells.Item(12, 18).Value = 1378.334817590824
$ws.Cells.Item(12, 19).Value = 0.01079816331466478
$ws.Cells.Item(12, 20).Value = 0.01171023614890298

$ws.Cells.Item(13, 7).Value = 57.98602933333333
$ws.Cells.Item(13, 8).Value = 173.958088
$ws.Cells.Item(13, 9).Value = 0.03631876156896331
$ws.Cells.Item(13, 10).Value = 0.03654368891224535
$ws.Cells.Item(13, 15).Value = 0.258118258924649
$ws.Cells.Item(13, 16).Value = 0.2781974433645584
$ws.Cells.Item(13, 17).Value = 132.9572683929573
$ws.Cells.Item(13, 18).Value = 1196.615415536616
$ws.Cells.Item(13, 19).Value = 0.009374535502480264
$ws.Cells.Item(13, 20).Value = 0.01016636082649642

$ws.Cells.Item(14, 7).Value = 57.98602933333333
$ws.Cells.Item(14, 8).Value = 173.958088
$ws.Cells.Item(14, 9).Value = 0.03631876156896331
$ws.Cells.Item(14, 10).Value = 0.03654368891224535
$ws.Cells.Item(14, 13).Value = 0.7354653333333333
$ws.Cells.Item(14, 14).Value = 2.206396
$ws.Cells.Item(14, 15).Value = 0.08279273334096697
$ws.Cells.Item(14, 16).Value = 0.08923323301721346
$ws.Cells.Item(14, 17).Value = 42.64671439231644
$ws.Cells.Item(14, 18).Value = 383.8204295308479
$ws.Cells.Item(14, 19).Value = 0.003006929541853339
$ws.Cells.Item(14, 20).Value = 0.003260911508014949

$ws.Cells.Item(15, 7).Value = 57.98602933333333
$ws.Cells.Item(15, 8).Value = 173.958088
$ws.Cells.Item(15, 9).Value = 0.03631876156896331
$ws.Cells.Item(15, 10).Value = 0.03654368891224535
$ws.Cells.Item(15, 13).Value = 1.9234645
$ws.Cells.Item(15, 14).Value = 3.846929
$ws.Cells.Item(15, 15).Value = 0.2165280621964276
$ws.Cells.Item(15, 16).Value = 0.1555812790893729
$ws.Cells.Item(15, 17).Value = 111.5340689186253
$ws.Cells.Item(15, 18).Value = 669.2044135117519
$ws.Cells.Item(15, 19).Value = 0.007864031063901713
$ws.Cells.Item(15, 20).Value = 0.005685513863611265

$ws.Cells.Item(16, 7).Value = 57.98602933333333
$ws.Cells.Item(16, 8).Value = 173.958088
$ws.Cells.Item(16, 9).Value = 0.03631876156896331
$ws.Cells.Item(16, 10).Value = 0.03654368891224535
$ws.Cells.Item(16, 13).Value = 1.290238
$ws.Cells.Item(16, 14).Value = 3.870714
$ws.Cells.Item(16, 15).Value = 0.1452445490479259
$ws.Cells.Item(16, 16).Value = 0.1565432154087437
$ws.Cells.Item(16, 17).Value = 74.81577851498133
$ws.Cells.Item(16, 18).Value = 673.3420066348319
$ws.Cells.Item(16, 19).Value = 0.005275102146063219
$ws.Cells.Item(16, 20).Value = 0.005720666565219741

$ws.Cells.Item(17, 7).Value = 29.481085
$ws.Cells.Item(17, 8).Value = 58.96217
$ws.Cells.Item(17, 9).Value = 0.01846507700595112
$ws.Cells.Item(17, 10).Value = 0.01238628926567028
$ws.Cells.Item(17, 13).Value = 2.641124333333333
$ws.Cells.Item(17, 14).Value = 7.923373
$ws.Cells.Item(17, 15).Value = 0.2973163964900306
$ws.Cells.Item(17, 16).Value = 0.3204448291201116
$ws.Cells.Item(17, 17).Value = 77.86321096656833
$ws.Cells.Item(17, 18).Value = 467.17926579941
$ws.Cells.Item(17, 19).Value = 0.005489970156320311
$ws.Cells.Item(17, 20).Value = 0.003969122347169984

$ws.Cells.Item(18, 7).Value = 29.481085
$ws.Cells.Item(18, 8).Value = 58.96217
$ws.Cells.Item(18, 9).Value = 0.01846507700595112
$ws.Cells.Item(18, 10).Value = 0.01238628926567028
$ws.Cells.Item(18, 15).Value = 0.258118258924649
$ws.Cells.Item(18, 16).Value = 0.2781974433645584
$ws.Cells.Item(18, 17).Value = 67.59773993711499
$ws.Cells.Item(18, 18).Value = 405.58643962269
$ws.Cells.Item(18, 19).Value = 0.004766173527685675
$ws.Cells.Item(18, 20).Value = 0.003445834006483345

$ws.Cells.Item(19, 7).Value = 29.481085
$ws.Cells.Item(19, 8).Value = 58.96217
$ws.Cells.Item(19, 9).Value = 0.01846507700595112
$ws.Cells.Item(19, 10).Value = 0.01238628926567028
$ws.Cells.Item(19, 13).Value = 0.7354653333333333
$ws.Cells.Item(19, 14).Value = 2.206396
$ws.Cells.Item(19, 15).Value = 0.08279273334096697
$ws.Cells.Item(19, 16).Value = 0.08923323301721346
$ws.Cells.Item(19, 17).Value = 21.68231600655333
$ws.Cells.Item(19, 18).Value = 130.09389603932
$ws.Cells.Item(19, 19).Value = 0.001528774196674132
$ws.Cells.Item(19, 20).Value = 0.001105268636262166

$ws.Cells.Item(20, 7).Value = 29.481085
$ws.Cells.Item(20, 8).Value = 58.96217
$ws.Cells.Item(20, 9).Value = 0.01846507700595112
$ws.Cells.Item(20, 10).Value = 0.01238628926567028
$ws.Cells.Item(20, 13).Value = 1.9234645
$ws.Cells.Item(20, 14).Value = 3.846929
$ws.Cells.Item(20, 15).Value = 0.2165280621964276
$ws.Cells.Item(20, 16).Value = 0.1555812790893729
$ws.Cells.Item(20, 17).Value = 56.7058204189825
$ws.Cells.Item(20, 18).Value = 226.82328167593
$ws.Cells.Item(20, 19).Value = 0.00399820734240641
$ws.Cells.Item(20, 20).Value = 0.001927074727123951

$ws.Cells.Item(21, 7).Value = 29.481085
$ws.Cells.Item(21, 8).Value = 58.96217
$ws.Cells.Item(21, 9).Value = 0.01846507700595112
$ws.Cells.Item(21, 10).Value = 0.01238628926567028
$ws.Cells.Item(21, 13).Value = 1.290238
$ws.Cells.Item(21, 14).Value = 3.870714
$ws.Cells.Item(21, 15).Value = 0.1452445490479259
$ws.Cells.Item(21, 16).Value = 0.1565432154087437
$ws.Cells.Item(21, 17).Value = 38.03761614823
$ws.Cells.Item(21, 18).Value = 228.22569688938
$ws.Cells.Item(21, 19).Value = 0.002681951782864597
$ws.Cells.Item(21, 20).Value = 0.001938989548630832

$ws.Cells.Item(22, 7).Value = 21.628479
$ws.Cells.Item(22, 8).Value = 64.885437
$ws.Cells.Item(22, 9).Value = 0.01354670393768061
$ws.Cells.Item(22, 10).Value = 0.01363060063446486
$ws.Cells.Item(22, 13).Value = 2.641124333333333
$ws.Cells.Item(22, 14).Value = 7.923373
$ws.Cells.Item(22, 15).Value = 0.2973163964900306
$ws.Cells.Item(22, 16).Value = 0.3204448291201116
$ws.Cells.Item(22, 17).Value = 57.12350217988899
$ws.Cells.Item(22, 18).Value = 514.111519619001
$ws.Cells.Item(22, 19).Value = 0.004027657199068507
$ws.Cells.Item(22, 20).Value = 0.004367855491115577

$ws.Cells.Item(23, 7).Value = 21.628479
$ws.Cells.Item(23, 8).Value = 64.885437
$ws.Cells.Item(23, 9).Value = 0.01354670393768061
$ws.Cells.Item(23, 10).Value = 0.01363060063446486
$ws.Cells.Item(23, 15).Value = 0.258118258924649
$ws.Cells.Item(23, 16).Value = 0.2781974433645584
$ws.Cells.Item(23, 17).Value = 49.592350440201
$ws.Cells.Item(23, 18).Value = 446.331153961809
$ws.Cells.Item(23, 19).Value = 0.003496651634561806
$ws.Cells.Item(23, 20).Value = 0.003791998248031453

$ws.Cells.Item(24, 7).Value = 21.628479
$ws.Cells.Item(24, 8).Value = 64.885437
$ws.Cells.Item(24, 9).Value = 0.01354670393768061
$ws.Cells.Item(24, 10).Value = 0.01363060063446486
$ws.Cells.Item(24, 13).Value = 0.7354653333333333
$ws.Cells.Item(24, 14).Value = 2.206396
$ws.Cells.Item(24, 15).Value = 0.08279273334096697
$ws.Cells.Item(24, 16).Value = 0.08923323301721346
$ws.Cells.Item(24, 17).Value = 15.906996517228
$ws.Cells.Item(24, 18).Value = 143.162968655052
$ws.Cells.Item(24, 19).Value = 0.001121568646761418
$ws.Cells.Item(24, 20).Value = 0.001216302562579781

$ws.Cells.Item(25, 7).Value = 21.628479
$ws.Cells.Item(25, 8).Value = 64.885437
$ws.Cells.Item(25, 9).Value = 0.01354670393768061
$ws.Cells.Item(25, 10).Value = 0.01363060063446486
$ws.Cells.Item(25, 13).Value = 1.9234645
$ws.Cells.Item(25, 14).Value = 3.846929
$ws.Cells.Item(25, 15).Value = 0.2165280621964276
$ws.Cells.Item(25, 16).Value = 0.1555812790893729
$ws.Cells.Item(25, 17).Value = 41.60161154549549
$ws.Cells.Item(25, 18).Value = 249.609669272973
$ws.Cells.Item(25, 19).Value = 0.002933241552774698
$ws.Cells.Item(25, 20).Value = 0.002120666281466461

$ws.Cells.Item(26, 7).Value = 21.628479
$ws.Cells.Item(26, 8).Value = 64.885437
$ws.Cells.Item(26, 9).Value = 0.01354670393768061
$ws.Cells.Item(26, 10).Value = 0.01363060063446486
$ws.Cells.Item(26, 13).Value = 1.290238
$ws.Cells.Item(26, 14).Value = 3.870714
$ws.Cells.Item(26, 15).Value = 0.1452445490479259
$ws.Cells.Item(26, 16).Value = 0.1565432154087437
$ws.Cells.Item(26, 17).Value = 27.905885488002
$ws.Cells.Item(26, 18).Value = 251.152969392018
$ws.Cells.Item(26, 19).Value = 0.001967584904514183
$ws.Cells.Item(26, 20).Value = 0.002133778051271591
